$p = $ppt.ActivePresentation
$c = $p.CustomXMLParts
Write-Output (Get-Member -InputObject $c)
